# Actualización de turnos en lista_de_turnos.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 ("Tránsito aéreo" shift): Francisca Ávila -> Francisca Zúñiga
$ws.Range("A14").Value = "Francisca Zúñiga"
$ws.Range("B14").Value = "18.392.207-6"

# Row 25 shift: Carlos Pinto -> Francisca Ávila
$ws.Range("A25").Value = "Francisca Ávila"
$ws.Range("B25").Value = "18.049.568-1"

# Move the active cell selection to G12 (matches the saved view state)
$ws.Range("G12").Select()
